$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D = 4 (Price), Column E = 5 (Volume(1h)).
# These columns store plain text (not numbers), including values
# that look numeric (e.g. "1.00", "301.60"). Excel auto-detects
# numeric-looking input as a real number, so we briefly force the
# "Text" number format while assigning the value, then restore the
# cell style to "Normal" so the saved file keeps the original
# (unstyled) look, exactly like the rest of the sheet.

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = "42.777.81"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(2, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.29%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.304.71"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(3, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.22%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(4, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.02%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "301.60"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(5, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.16%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "96.21"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(6, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.24%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(7, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.37%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(8, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.05%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(9, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.11%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(10, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.10%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "19.28"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(11, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.07%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0790"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(12, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.03%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(13, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.08%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.78"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(14, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.89%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.661.31"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(15, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.00%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.299.94"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(16, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.63%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(17, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.21%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "42.709.23"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(18, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.34%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "12.27"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(19, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -5.71%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0893"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(20, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.39%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "67.81"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(22, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +1.04%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(23, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +6.77%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "235.32"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(24, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.31%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(25, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.15%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.42"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(26, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.43%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "24.42"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(27, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.73%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.37"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(28, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +15.00%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "165.17"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(29, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.56%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(30, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.48%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(31, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.24%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(32, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.02%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.99"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(33, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.28%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "17.60"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(34, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.25%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.48"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(35, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -5.82%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0704"
$cell.Style = "Normal"

$cell = $ws.Cells.Item(37, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -2.95%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(38, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.73%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(39, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.75%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(40, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.60%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.109"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(41, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.75%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "20.43"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(42, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +13.16%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.976.23"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(43, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.02%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(44, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +5.16%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(45, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.18%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(46, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.90%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(47, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.54%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.529.75"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(48, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.93%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "53.40"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(49, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -0.32%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.79"
$cell.Style = "Normal"
$cell = $ws.Cells.Item(50, 5)
$cell.NumberFormat = "@"
$cell.Value = "  -1.90%  "
$cell.Style = "Normal"

$cell = $ws.Cells.Item(51, 5)
$cell.NumberFormat = "@"
$cell.Value = "  +0.34%  "
$cell.Style = "Normal"
